$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename column headers (row 1, B1:M1) ---
$headers = @{
    "B1" = "CHCPIY_DIFF"
    "C1" = "CHJOB_DIFF"
    "D1" = "CHPMI_LOG"
    "E1" = "CNCPI_DIFF"
    "F1" = "CNPMIB_LOG"
    "G1" = "EUHICY_DIFF"
    "H1" = "EUUNR_DIFF"
    "I1" = "RUCPIY_DIFF"
    "J1" = "RUUNR_DIFF"
    "K1" = "USCPI_DIFF"
    "L1" = "USPMI_LOG"
    "M1" = "USUNR_DIFF"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- Rename row labels (column A, A2:A13), same names in row order ---
$rowLabels = @(
    "CHCPIY_DIFF",
    "CHJOB_DIFF",
    "CHPMI_LOG",
    "CNCPI_DIFF",
    "CNPMIB_LOG",
    "EUHICY_DIFF",
    "EUUNR_DIFF",
    "RUCPIY_DIFF",
    "RUUNR_DIFF",
    "USCPI_DIFF",
    "USPMI_LOG",
    "USUNR_DIFF"
)
for ($i = 0; $i -lt $rowLabels.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $rowLabels[$i]
}

# --- Update data cell values (B2:M13) ---
$data = @{
    2 = @{ "C" = 2.12; "D" = 1.19; "E" = 1.18; "F" = 1.31; "G" = 1.3; "H" = 1.27; "I" = 1.06; "J" = 1.42; "K" = 1.58; "L" = 1.66; "M" = 2.4 }
    3 = @{ "B" = 1.51; "D" = 1.21; "E" = 1.15; "F" = 1.09; "G" = 1.4; "H" = 1.27; "I" = 1.07; "J" = 1.47; "K" = 1.6; "L" = 1.65; "M" = 1.86 }
    4 = @{ "B" = 1.6; "C" = 2.27; "E" = 1.18; "F" = 1.27; "G" = 1.42; "H" = 1.29; "I" = 1.06; "J" = 1.46; "K" = 1.6; "L" = 1.68; "M" = 2.35 }
    5 = @{ "B" = 1.62; "C" = 2.21; "D" = 1.2; "F" = 1.33; "G" = 1.41; "H" = 1.3; "I" = 1.05; "J" = 1.46; "K" = 1.58; "L" = 1.63; "M" = 2.37 }
    6 = @{ "B" = 1.59; "C" = 1.86; "D" = 1.15; "E" = 1.18; "G" = 1.42; "H" = 1.28; "I" = 1.08; "J" = 1.46; "K" = 1.62; "L" = 1.68; "M" = 2.3 }
    7 = @{ "B" = 1.49; "C" = 2.25; "D" = 1.21; "E" = 1.18; "F" = 1.33; "H" = 1.29; "I" = 1.08; "J" = 1.47; "K" = 1.5; "L" = 1.68; "M" = 2.4 }
    8 = @{ "B" = 1.58; "C" = 2.22; "D" = 1.19; "E" = 1.18; "F" = 1.31; "G" = 1.41; "I" = 1.08; "J" = 1.37; "K" = 1.57; "L" = 1.61; "M" = 2.41 }
    9 = @{ "B" = 1.6; "C" = 2.26; "D" = 1.18; "E" = 1.16; "F" = 1.33; "G" = 1.42; "H" = 1.3; "J" = 1.47; "K" = 1.62; "L" = 1.68; "M" = 2.42 }
    10 = @{ "B" = 1.57; "C" = 2.27; "D" = 1.2; "E" = 1.18; "F" = 1.33; "G" = 1.42; "H" = 1.22; "I" = 1.08; "K" = 1.62; "L" = 1.65; "M" = 2.11 }
    11 = @{ "B" = 1.58; "C" = 2.24; "D" = 1.19; "E" = 1.16; "F" = 1.33; "G" = 1.31; "H" = 1.27; "I" = 1.08; "J" = 1.47; "L" = 1.64; "M" = 2.41 }
    12 = @{ "B" = 1.59; "C" = 2.23; "D" = 1.2; "E" = 1.14; "F" = 1.33; "G" = 1.42; "H" = 1.25; "I" = 1.08; "J" = 1.44; "K" = 1.58; "M" = 1.83 }
    13 = @{ "B" = 1.6; "C" = 1.75; "D" = 1.17; "E" = 1.16; "F" = 1.27; "G" = 1.41; "H" = 1.3; "I" = 1.08; "J" = 1.28; "K" = 1.62; "L" = 1.28 }
}
foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $data[$row][$col]
    }
}
